$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 648.25
$ws.Range("J42").Value = 533
$ws.Range("L42").Value = 1599
$ws.Range("N42").Value = -2059
$ws.Range("H47").Value = 66499.5
$ws.Range("J47").Value = 66499.5
$ws.Range("L47").Value = 66499.5
$ws.Range("N47").Value = -68443.5
$ws.Range("H74").Value = 8125.125
$ws.Range("I74").Value = 6250.25
$ws.Range("K74").Value = 6250.25
$ws.Range("M74").Value = -5314.25
$ws.Range("H77").Value = 8125.125
$ws.Range("I77").Value = 6250.25
$ws.Range("K77").Value = 31251.25
$ws.Range("M77").Value = -26571.25
$ws.Range("H98").Value = 12017.846
$ws.Range("I98").Value = 13838.728
$ws.Range("J98").Value = 2003
$ws.Range("K98").Value = 13838.728
$ws.Range("L98").Value = 2003
$ws.Range("M98").Value = -12340.728
$ws.Range("N98").Value = -4999
$ws.Range("H122").Value = 12017.846
$ws.Range("I122").Value = 13838.728
$ws.Range("J122").Value = 2003
$ws.Range("K122").Value = 41516.18399999999
$ws.Range("L122").Value = 6009
$ws.Range("M122").Value = -39066.18399999999
$ws.Range("N122").Value = -10909
$ws.Range("H137").Value = 7581489
$ws.Range("I137").Value = 9618451
$ws.Range("K137").Value = 28855353
$ws.Range("M137").Value = -28852803
$ws.Range("H138").Value = 3538.0576
$ws.Range("I138").Value = 2131.7273
$ws.Range("K138").Value = 6395.1819
$ws.Range("M138").Value = -1255.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2557065.2
$ws.Range("I2").Value = 3286809
$ws.Range("K2").Value = 3286809
$ws.Range("M2").Value = -3286696
$ws.Range("H28").Value = 18967.143
$ws.Range("I28").Value = 18967.143
$ws.Range("K28").Value = 18967.143
$ws.Range("M28").Value = -18775.143
$ws.Range("H32").Value = 2918248.5
$ws.Range("I32").Value = 3145598
$ws.Range("K32").Value = 3145598
$ws.Range("M32").Value = -3145311
$ws.Range("H37").Value = 36246.75
$ws.Range("I37").Value = 14999.667
$ws.Range("K37").Value = 14999.667
$ws.Range("M37").Value = -14726.667
$ws.Range("H50").Value = 4363.8
$ws.Range("I50").Value = 5341
$ws.Range("J50").Value = 2898
$ws.Range("K50").Value = 5341
$ws.Range("L50").Value = 2898
$ws.Range("M50").Value = -4627
$ws.Range("N50").Value = -4326
$ws.Range("H74").Value = 211050.05
$ws.Range("I74").Value = 239176.36
$ws.Range("K74").Value = 239176.36
$ws.Range("M74").Value = -238302.36
$ws.Range("H77").Value = 211050.05
$ws.Range("I77").Value = 239176.36
$ws.Range("K77").Value = 1195881.8
$ws.Range("M77").Value = -1191513.8
$ws.Range("H99").Value = 18967.143
$ws.Range("I99").Value = 18967.143
$ws.Range("K99").Value = 18967.143
$ws.Range("M99").Value = -15972.143
$ws.Range("H116").Value = 2557065.2
$ws.Range("I116").Value = 3286809
$ws.Range("K116").Value = 3286809
$ws.Range("M116").Value = -3284515

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2557065.2
$ws.Range("I3").Value = 3286809
$ws.Range("K3").Value = 3286809
$ws.Range("M3").Value = -3286695
$ws.Range("H86").Value = 182997.81
$ws.Range("I86").Value = 1052.8889
$ws.Range("K86").Value = 1052.8889
$ws.Range("M86").Value = 70.11110000000008
$ws.Range("H89").Value = 182997.81
$ws.Range("I89").Value = 1052.8889
$ws.Range("K89").Value = 5264.4445
$ws.Range("M89").Value = 351.5555000000004
$ws.Range("H107").Value = 1072.6428
$ws.Range("I107").Value = 988.28
$ws.Range("K107").Value = 988.28
$ws.Range("M107").Value = 931.72
$ws.Range("H134").Value = 3779.1614
$ws.Range("I134").Value = 2440.24
$ws.Range("K134").Value = 7320.719999999999
$ws.Range("M134").Value = -4785.719999999999
$ws.Range("H135").Value = 88389.5
$ws.Range("J135").Value = 88389.5
$ws.Range("L135").Value = 88389.5
$ws.Range("N135").Value = -98529.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 557.1515000000001
$ws.Range("I19").Value = 298
$ws.Range("K19").Value = 298
$ws.Range("M19").Value = -128
$ws.Range("H24").Value = 557.1515000000001
$ws.Range("I24").Value = 298
$ws.Range("K24").Value = 298
$ws.Range("M24").Value = -128
$ws.Range("H52").Value = 15000
$ws.Range("I52").Value = 15000
$ws.Range("K52").Value = 15000
$ws.Range("M52").Value = -14706
$ws.Range("H93").Value = 10030.2
$ws.Range("J93").Value = 39991
$ws.Range("L93").Value = 39991
$ws.Range("N93").Value = -43735
$ws.Range("H104").Value = 24165.666
$ws.Range("I104").Value = 6250
$ws.Range("K104").Value = 6250
$ws.Range("M104").Value = -3629
$ws.Range("H115").Value = 94979
$ws.Range("J115").Value = 94979
$ws.Range("L115").Value = 94979
$ws.Range("N115").Value = -97329
$ws.Range("H132").Value = 2462.7407
$ws.Range("I132").Value = 2507.5386
$ws.Range("K132").Value = 7522.6158
$ws.Range("M132").Value = -4992.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2396.0557
$ws.Range("J34").Value = 3066.3572
$ws.Range("L34").Value = 9199.071599999999
$ws.Range("N34").Value = -9367.071599999999
$ws.Range("H39").Value = 7335.3335
$ws.Range("H55").Value = 87.5
$ws.Range("I55").Value = 87.5
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 262.5
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -85.5
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 51535.57
$ws.Range("J32").Value = 70499.75
$ws.Range("L32").Value = 70499.75
$ws.Range("N32").Value = -71091.75
$ws.Range("H42").Value = 76145
$ws.Range("J42").Value = 76145
$ws.Range("L42").Value = 76145
$ws.Range("N42").Value = -77115
$ws.Range("H70").Value = 7199.5713
$ws.Range("I70").Value = 6649.25
$ws.Range("K70").Value = 6649.25
$ws.Range("M70").Value = -6379.25
$ws.Range("H73").Value = 7199.5713
$ws.Range("I73").Value = 6649.25
$ws.Range("K73").Value = 6649.25
$ws.Range("M73").Value = -5713.25
$ws.Range("H114").Value = 79910
$ws.Range("J114").Value = 79910
$ws.Range("L114").Value = 79910
$ws.Range("N114").Value = -88588
$ws.Range("H115").Value = 76145
$ws.Range("J115").Value = 76145
$ws.Range("L115").Value = 76145
$ws.Range("N115").Value = -78495
$ws.Range("H122").Value = 4769.75
$ws.Range("I122").Value = 2682.1177
$ws.Range("J122").Value = 9839.714
$ws.Range("K122").Value = 8046.353099999999
$ws.Range("L122").Value = 29519.142
$ws.Range("M122").Value = -5596.353099999999
$ws.Range("N122").Value = -34419.142
$ws.Range("H132").Value = 4273.92
$ws.Range("I132").Value = 2867.5881
$ws.Range("J132").Value = 7262.375
$ws.Range("K132").Value = 8602.764299999999
$ws.Range("L132").Value = 21787.125
$ws.Range("M132").Value = -6072.764299999999
$ws.Range("N132").Value = -26847.125
$ws.Range("H137").Value = 71500
$ws.Range("J137").Value = 71500
$ws.Range("L137").Value = 71500
$ws.Range("N137").Value = -81700

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 204900.6
$ws.Range("I7").Value = 337067.66
$ws.Range("K7").Value = 337067.66
$ws.Range("M7").Value = -336955.66
$ws.Range("H16").Value = 3016.12
$ws.Range("I16").Value = 2433.7368
$ws.Range("J16").Value = 4860.3335
$ws.Range("K16").Value = 2433.7368
$ws.Range("L16").Value = 4860.3335
$ws.Range("M16").Value = -2263.7368
$ws.Range("N16").Value = -5200.3335
$ws.Range("H22").Value = 3435.2285
$ws.Range("I22").Value = 2006.7059
$ws.Range("J22").Value = 4784.3887
$ws.Range("K22").Value = 2006.7059
$ws.Range("L22").Value = 4784.3887
$ws.Range("M22").Value = -1711.7059
$ws.Range("N22").Value = -5374.3887
$ws.Range("H27").Value = 3435.2285
$ws.Range("I27").Value = 2006.7059
$ws.Range("J27").Value = 4784.3887
$ws.Range("K27").Value = 2006.7059
$ws.Range("L27").Value = 4784.3887
$ws.Range("M27").Value = -1899.7059
$ws.Range("N27").Value = -4998.3887
$ws.Range("H126").Value = 204900.6
$ws.Range("I126").Value = 337067.66
$ws.Range("K126").Value = 1011202.98
$ws.Range("M126").Value = -1008732.98
$ws.Range("H132").Value = 3202.9285
$ws.Range("I132").Value = 1550.862
$ws.Range("K132").Value = 4652.586
$ws.Range("M132").Value = -2122.586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 7553
$ws.Range("I11").Value = 1999
$ws.Range("J11").Value = 8346.429
$ws.Range("K11").Value = 1999
$ws.Range("L11").Value = 8346.429
$ws.Range("M11").Value = -1857
$ws.Range("N11").Value = -8630.429
$ws.Range("H122").Value = 1678.72
$ws.Range("I122").Value = 1520.3478
$ws.Range("K122").Value = 4561.0434
$ws.Range("M122").Value = -2111.0434
$ws.Range("H127").Value = 33333
$ws.Range("I127").Value = 33333
$ws.Range("K127").Value = 33333
$ws.Range("M127").Value = -28373
$ws.Range("H132").Value = 3747.7354
$ws.Range("I132").Value = 2011.5714
$ws.Range("K132").Value = 6034.7142
$ws.Range("M132").Value = -3504.7142
$ws.Range("H136").Value = 2000.5555
$ws.Range("I136").Value = 809.8182
$ws.Range("K136").Value = 2429.4546
$ws.Range("M136").Value = 120.5454
